$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1.55
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1.06
$ws.Range("H2").Value = 0.9
$ws.Range("I2").Value = 0.7

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1.62
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0.43

# Row 4
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 3
$ws.Range("I4").Value = 0.77
